$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "ilz ont affayre. Il sont par dehors touts d" ->
#         "ilz ont affayre. Il s" + <del>e</del> + "ont par dehors touts d"
# The visible letters are unchanged; we only splice a red, small
# Courier-New "<del>" / "</del>" markup pair around a plain "e" right
# between "Il s" and "ont par dehors touts d".
# ---------------------------------------------------------------------------

$anchor = $d.Content
$found = $anchor.Find.Execute("ilz ont affayre. Il s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find insertion anchor for edit 1"
}
$anchor.Collapse(0)
$pos = $anchor.Start

# Insert the whole markup as plain text first (it will merge into the
# surrounding run because the formatting is identical at this point), then
# go back and re-colour just the "<del>" and "</del>" tag pieces, leaving
# the "e" itself with the plain, inherited formatting.
$anchor.InsertAfter("<del>e</del>")

$delOpenStart = $pos
$delOpenEnd = $pos + 5          # "<del>"
$delCloseStart = $pos + 6       # after the "e"
$delCloseEnd = $pos + 12        # "</del>"

$rOpen = $d.Range($delOpenStart, $delOpenEnd)
$rOpen.Font.Name = "Courier New"
$rOpen.Font.Color = 1118633
$rOpen.Font.Size = 9

$rClose = $d.Range($delCloseStart, $delCloseEnd)
$rClose.Font.Name = "Courier New"
$rClose.Font.Color = 1118633
$rClose.Font.Size = 9

# ---------------------------------------------------------------------------
# Edit 2: "ilz ne se crevent poinct." -> "ilz ne se crevent poinct. "
# (trailing space appended, same run formatting)
# ---------------------------------------------------------------------------

$tail = $d.Content
$found2 = $tail.Find.Execute("ilz ne se crevent poinct.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target text for edit 2"
}
$tail.Collapse(0)
$tail.InsertAfter(" ")

Write-Output "done"
